$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 - 四方坪站
$ws.Cells.Item(56, 1).Value = 45958
$ws.Cells.Item(56, 2).Value = "四方坪站"
$ws.Cells.Item(56, 3).Value = 8548.61
$ws.Cells.Item(56, 4).Value = 6884.26
$ws.Cells.Item(56, 5).Value = 2855.65
$ws.Cells.Item(56, 6).Value = 383

# Row 57 - 高岭站
$ws.Cells.Item(57, 1).Value = 45958
$ws.Cells.Item(57, 2).Value = "高岭站"
$ws.Cells.Item(57, 3).Value = 4029.96
$ws.Cells.Item(57, 4).Value = 3331.48
$ws.Cells.Item(57, 5).Value = 995.91
$ws.Cells.Item(57, 6).Value = 142

# Move the active selection to match the edit (I56)
$ws.Range("I56").Select()
